$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 21 (pushes current rows 21-28 down to 22-29),
# inheriting formatting from the row above.
$ws.Rows.Item(21).Insert()

# Resize the DataProvider table (and its AutoFilter) to include the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B7:G29"))

# Populate the new row 21 with the new "127.0.0.1" IPv4 failure test case.
$ws.Range("B21").Value = "{PREVIOUS}"
$ws.Range("E21").Value = "127.0.0.1"
$ws.Range("F21").Value = "fail"
$ws.Range("G21").Value = "EPP_HOST_CREATE_SERVER_ACCEPTS_INVALID_IPV4_ADDRESS"

# "true" is a reserved boolean-like literal, so a plain .Value assignment
# gets auto-coerced to a Boolean cell. Force it in as text (quote-prefix),
# then repaste the formatting from the row above to drop the quote-prefix
# style while keeping the cell typed as text.
$ws.Range("C21").Value = "'true"
$ws.Range("C20").Copy()
$ws.Range("C21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the selection/view to match the saved state.
$ws.Range("E21").Select() | Out-Null
$excel.ActiveWindow.Zoom = 130
